$d = $word.ActiveDocument

# The document contains five occurrences of an <id> tag whose value and
# closing tag were previously split across separate runs (different
# character formatting), e.g.:
#   <id>  (Courier New, color 7f6000)  +  p098v_1  (plain, color 000000)  +  </id>  (Courier New, color 7f6000)
# They should be merged into a single run/text node reading "<id>p098v_1</id>".
# Doing a MatchCase, whole-text Find & Replace across the three runs makes
# Word collapse them into one run using the formatting of the first
# matched character (the opening "<id>" run's Courier New formatting).
foreach ($n in 1..5) {
    $text = "<id>p098v_$n</id>"
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $text, 2)
}

Write-Output "done"
